$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Tree")

# 1. Delete column B ("Source") - everything shifts one column left (C->B, D->C, ... J->I)
$ws.Columns("B").Delete()

# 2. Fix the pre-existing hyperlink whose range metadata still points at the old
#    (pre-shift) address D38 even though the cell itself is now C38.
$ws.Range("D38").Hyperlinks.Delete()
$leetUrl = "https://leetcode.com/problems/all-nodes-distance-k-in-binary-tree/"
$ws.Hyperlinks.Add($ws.Range("C38"), $leetUrl, "", "", $leetUrl)
$ws.Range("C38").Value = "All Nodes Distance K in Binary Tree - LeetCode"
$ws.Range("C38").Style = "Hyperlink"

# 3. Add new "Solution link" hyperlinks in column I (rows 4-9), in the same
#    order the original author inserted them.
$base = "https://github.com/spartan4cs/CP/blob/main/2.Pepcoding/TSP2/Level1/7.BinaryTree/"

$u = $base + "BinaryTreeConstruction.java"
$ws.Hyperlinks.Add($ws.Range("I4"), $u, "", "", $u)
$ws.Range("I4").Value = "CP/BinaryTreeConstruction.java at main " + [char]0x00B7 + " spartan4cs/CP (github.com)"
$ws.Range("I4").Style = "Hyperlink"

$u = $base + "BinaryTreeConstruction.java"
$ws.Hyperlinks.Add($ws.Range("I6"), $u, "", "", $u)
$ws.Range("I6").Value = "CP/BinaryTreeConstruction.java at main " + [char]0x00B7 + " spartan4cs/CP (github.com)"
$ws.Range("I6").Style = "Hyperlink"

$u = $base + "BinaryTreeConstruction.java"
$ws.Hyperlinks.Add($ws.Range("I5"), $u, "", "", $u)
$ws.Range("I5").Value = "CP/BinaryTreeConstruction.java at main " + [char]0x00B7 + " spartan4cs/CP (github.com)"
$ws.Range("I5").Style = "Hyperlink"

$u = $base + "SizeMinMaxHgt.java"
$ws.Hyperlinks.Add($ws.Range("I7"), $u, "", "", $u)
$ws.Range("I7").Value = "CP/SizeMinMaxHgt.java at main " + [char]0x00B7 + " spartan4cs/CP (github.com)"
$ws.Range("I7").Style = "Hyperlink"

$u = $base + "LevelOrderTraversal.java"
$ws.Hyperlinks.Add($ws.Range("I9"), $u, "", "", $u)
$ws.Range("I9").Value = "CP/LevelOrderTraversal.java at main " + [char]0x00B7 + " spartan4cs/CP (github.com)"
$ws.Range("I9").Style = "Hyperlink"

$u = $base + "Traversal.java"
$ws.Hyperlinks.Add($ws.Range("I8"), $u, "", "", $u)
$ws.Range("I8").Value = "CP/Traversal.java at main " + [char]0x00B7 + " spartan4cs/CP (github.com)"
$ws.Range("I8").Style = "Hyperlink"

# 4. Row 9 got manually resized (taller) by the author when they added the
#    new column to accommodate wrapped text.
$ws.Rows(9).RowHeight = 30

# 5. Shrink the new "Solution link" column width.
$ws.Columns("I").ColumnWidth = 13.6

# 6. The selection cursor ends on C10 in the saved file.
$ws.Range("C10").Select()
